$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F width adjustment
$ws.Range("F1").ColumnWidth = 12.484375

# zone column (F) updates: AUTRE -> BARGNY
$ws.Range("F73:F114").Value = "BARGNY"
$ws.Range("F126:F173").Value = "BARGNY"

# zone column (F) updates: AUTRE -> MTOA
$ws.Range("F239:F289").Value = "MTOA"
$ws.Range("F301:F337").Value = "MTOA"

# zone column (F) updates: MEDINA -> RUFISQUE, AUTRE -> RUFISQUE
$ws.Range("F490:F535").Value = "RUFISQUE"
$ws.Range("F536:F578").Value = "RUFISQUE"

# Row 120 quantity/amount correction
$ws.Range("J120").Value = 5.0
$ws.Range("K120").Value = 61250.0
